$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text, matching the source data
# (values like "30.571.68" are not valid numbers, and values like "1.011"
# must not be reinterpreted/rounded as floating point numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.571.68"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.111.88"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.95"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5255"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4577"
$ws.Range("E8").Value = "  +5.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.14"
$ws.Range("E9").Value = "  +5.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08960"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.174"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.42"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.120.31"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.868"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.195"
$ws.Range("E15").Value = "  +6.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001177"
$ws.Range("E16").Value = "  +5.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.01"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.012"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06681"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.22"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.244"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.637.06"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.73"
$ws.Range("E24").Value = "  +4.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.363"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.360.54"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.33"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.47"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.531"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.69"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.631"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.323"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.963"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.46"
$ws.Range("E36").Value = "  +3.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.898"
$ws.Range("E37").Value = "  +8.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02591"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06838"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2312"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.58"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6866"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.255"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.99"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.681"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.249"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000347"
$ws.Range("E49").Value = "  +14.12%  "

# Row 31: only E changes
$ws.Range("E31").Value = "  +2.19%  "

# Row 44/45: NEARProtocol and Decentraland swap positions, with updated data
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6454"
$ws.Range("E44").Value = "  +1.28%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.323"
$ws.Range("E45").Value = "  +5.37%  "

# Row 50: only E changes
$ws.Range("E50").Value = "  -2.08%  "

# Row 51: Aave -> WOONetwork
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3389"
$ws.Range("E51").Value = "  +12.92%  "
